$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$updates_ALC = @{
    "H5" = 109.083336
    "I5" = 106
    "J5" = 118.333336
    "K5" = 106
    "L5" = 118.333336
    "M5" = 9
    "N5" = -348.333336
    "H43" = 5300.1665
    "I43" = 4650.75
    "K43" = 4650.75
    "M43" = -4581.75
    "H58" = 4293.857
    "I58" = 302
    "J58" = 6511.5557
    "K58" = 906
    "L58" = 19534.6671
    "M58" = -756
    "N58" = -19834.6671
    "H62" = 3375.647
    "I62" = 3292.6
    "K62" = 3292.6
    "M62" = -2668.6
    "H65" = 3375.647
    "I65" = 3292.6
    "K65" = 16463
    "M65" = -13343
    "H69" = 17036.715
    "I69" = 14504.333
    "J69" = 17727.363
    "K69" = 43512.999
    "L69" = 53182.08900000001
    "M69" = -42638.999
    "N69" = -54930.08900000001
    "H72" = 17036.715
    "I72" = 14504.333
    "J72" = 17727.363
    "K72" = 130538.997
    "L72" = 159546.267
    "M72" = -126170.997
    "N72" = -168282.267
    "H82" = 549
    "I82" = 549
    "K82" = 1647
    "M82" = -1241
    "H85" = 549
    "I85" = 549
    "K85" = 1647
    "M85" = -243
    "H137" = 15874939
    "J137" = 2320.25
    "L137" = 6960.75
    "N137" = -12060.75
}
foreach ($addr in $updates_ALC.Keys) {
    $ws.Range($addr).Value = $updates_ALC[$addr]
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$updates_ARM = @{
    "H32" = 15055.826
    "I32" = 15814.45
    "J32" = 9998.333000000001
    "K32" = 15814.45
    "L32" = 9998.333000000001
    "M32" = -15527.45
    "N32" = -10572.333
    "H74" = 3373.4443
    "J74" = 4402.8
    "L74" = 4402.8
    "N74" = -6150.8
    "H77" = 3373.4443
    "J77" = 4402.8
    "L77" = 22014
    "N77" = -30750
    "H132" = 3614.4167
    "I132" = 2584.516
    "J132" = 9999.799999999999
    "K132" = 7753.548000000001
    "L132" = 29999.4
    "M132" = -5223.548000000001
    "N132" = -35059.39999999999
    "H122" = 3018.75
    "I122" = 3018.75
    "J122" = 0
    "K122" = 9056.25
    "L122" = 0
    "M122" = -6606.25
}
foreach ($addr in $updates_ARM.Keys) {
    $ws.Range($addr).Value = $updates_ARM[$addr]
}
$clears_ARM = @("N122")
foreach ($addr in $clears_ARM) {
    $ws.Range($addr).ClearContents()
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$updates_BSM = @{
    "H22" = 875.625
    "J22" = 747.5
    "L22" = 747.5
    "N22" = -1093.5
    "H86" = 40428
    "I86" = 43088.625
    "J86" = 33333
    "K86" = 43088.625
    "L86" = 33333
    "M86" = -41965.625
    "N86" = -35579
    "H89" = 40428
    "I89" = 43088.625
    "J89" = 33333
    "K89" = 215443.125
    "L89" = 166665
    "M89" = -209827.125
    "N89" = -177897
    "H134" = 1784.5883
    "I134" = 1771.125
    "K134" = 5313.375
    "M134" = -2778.375
}
foreach ($addr in $updates_BSM.Keys) {
    $ws.Range($addr).Value = $updates_BSM[$addr]
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$updates_CRP = @{
    "H31" = 7973.8335
    "I31" = 14635.875
    "J31" = 4642.8125
    "K31" = 14635.875
    "L31" = 4642.8125
    "M31" = -14340.875
    "N31" = -5232.8125
    "H34" = 7973.8335
    "I34" = 14635.875
    "J34" = 4642.8125
    "K34" = 14635.875
    "L34" = 4642.8125
    "M34" = -14433.875
    "N34" = -5046.8125
    "H96" = 20849.2
    "J96" = 20849.2
    "L96" = 20849.2
    "N96" = -26341.2
    "H134" = 1642.409
    "I134" = 1747.3889
    "J134" = 1170
    "K134" = 5242.1667
    "L134" = 3510
    "M134" = -2707.1667
    "N134" = -8580
    "H54" = 9000
    "I54" = 9000
    "K54" = 9000
    "M54" = -8342
}
foreach ($addr in $updates_CRP.Keys) {
    $ws.Range($addr).Value = $updates_CRP[$addr]
}

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$updates_CUL = @{
    "H69" = 1887
    "I69" = 1223.25
    "J69" = 2108.25
    "K69" = 3669.75
    "L69" = 6324.75
    "M69" = -2858.75
    "N69" = -7946.75
    "H72" = 1887
    "I72" = 1223.25
    "J72" = 2108.25
    "K72" = 11009.25
    "L72" = 18974.25
    "M72" = -6953.25
    "N72" = -27086.25
    "H95" = 100000
    "J95" = 100000
    "L95" = 300000
    "N95" = -304118
    "H113" = 2068.3809
    "I113" = 2238.3
    "J113" = 1913.909
    "K113" = 6714.900000000001
    "L113" = 5741.727000000001
    "M113" = -4544.900000000001
    "N113" = -10081.727
    "H74" = 5000
    "J74" = 5000
    "L74" = 15000
    "H77" = 5000
    "J77" = 5000
    "L77" = 45000
    "N74" = -17122
    "N77" = -55608
}
foreach ($addr in $updates_CUL.Keys) {
    $ws.Range($addr).Value = $updates_CUL[$addr]
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$updates_GSM = @{
    "H132" = 9261741
    "I132" = 2211.9656
    "J132" = 47622650
    "K132" = 6635.8968
    "L132" = 142867950
    "M132" = -4105.8968
    "N132" = -142873010
    "H95" = 8000
    "J95" = 8000
    "L95" = 8000
    "H80" = 2500
    "I80" = 2500
    "J80" = 0
    "K80" = 2500
    "L80" = 0
    "M80" = -1502
    "H83" = 2500
    "I83" = 2500
    "J83" = 0
    "K83" = 12500
    "L83" = 0
    "M83" = -7508
    "H122" = 200
    "I122" = 200
    "J122" = 0
    "K122" = 600
    "L122" = 0
    "M122" = 1850
    "N95" = -13492
}
foreach ($addr in $updates_GSM.Keys) {
    $ws.Range($addr).Value = $updates_GSM[$addr]
}
$clears_GSM = @("N80", "N83", "N122")
foreach ($addr in $clears_GSM) {
    $ws.Range($addr).ClearContents()
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$updates_LTW = @{
    "H22" = 100003940
    "J22" = 250002450
    "L22" = 250002450
    "N22" = -250003040
    "H27" = 100003940
    "J27" = 250002450
    "L27" = 250002450
    "N27" = -250002664
    "H46" = 5480.067
    "I46" = 2966.25
    "J46" = 6394.1816
    "K46" = 2966.25
    "L46" = 6394.1816
    "M46" = -2778.25
    "N46" = -6770.1816
    "H132" = 2946.4
    "I132" = 2974.5
    "J132" = 2914.2856
    "K132" = 8923.5
    "L132" = 8742.856800000001
    "M132" = -6393.5
    "N132" = -13802.8568
    "H136" = 166668300
    "I136" = 2450
    "J136" = 500000000
    "K136" = 7350
    "L136" = 1500000000
    "M136" = -4800
    "N136" = -1500005100
    "H122" = 2834.923
    "I122" = 2834.923
    "J122" = 0
    "K122" = 8504.769
    "L122" = 0
    "M122" = -6054.769
}
foreach ($addr in $updates_LTW.Keys) {
    $ws.Range($addr).Value = $updates_LTW[$addr]
}
$clears_LTW = @("N122")
foreach ($addr in $clears_LTW) {
    $ws.Range($addr).ClearContents()
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$updates_WVR = @{
    "H122" = 2159.75
    "I122" = 2159.75
    "K122" = 6479.25
    "M122" = -4029.25
    "H125" = 25665.334
    "J125" = 25665.334
    "L125" = 25665.334
    "N125" = -35505.334
    "H132" = 166668670
    "I132" = 2349.5
    "J132" = 500001300
    "K132" = 7048.5
    "L132" = 1500003900
    "M132" = -4518.5
    "N132" = -1500008960
    "H136" = 8847.929
    "I136" = 8847.929
    "K136" = 26543.787
    "M136" = -23993.787
}
foreach ($addr in $updates_WVR.Keys) {
    $ws.Range($addr).Value = $updates_WVR[$addr]
}
